# Apply the resume contact-info update:
#   - phone number: 779-7473 -> 262-5832
#   - email: bfs5030@gmail.com -> byron@hackbytes.com
#   - move the "_GoBack" bookmark from the Applications line to sit
#     immediately before the new email run

$d = $word.ActiveDocument

# 1. Update phone number (personal -> business number)
$d.Content.Find.Execute("779-7473", $true, $false, $false, $false, $false, $true, 1, $false, "262-5832", 2) | Out-Null

# 2. Update email address
$d.Content.Find.Execute("bfs5030@gmail.com", $true, $false, $false, $false, $false, $true, 1, $false, "byron@hackbytes.com", 2) | Out-Null

# 3. Remove the existing hidden "_GoBack" bookmark (currently sitting after
#    "Xcode, Eclipse, Fossil, Git") -- Word relocates this automatically to
#    mark the last edit position, so we delete the stale one...
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 4. ...and re-add it immediately before the newly-typed email text, which is
#    where the author's last edit actually happened.
$emailRange = $d.Content
$emailRange.Find.Execute("byron@hackbytes.com") | Out-Null
$bookmarkPos = $d.Range($emailRange.Start, $emailRange.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkPos) | Out-Null
